$d = $word.ActiveDocument

# The document currently ends with a single trailing empty paragraph.
# Turn that empty paragraph into the first new paragraph of text, then
# keep appending the rest of the new content after it.
$p = $d.Paragraphs.Last
$p.Range.InsertAfter("My final years of graduate school, I started writing this list that I wanted to embody my experience.")

# "PhD realizations:" intro paragraph
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.InsertAfter("PhD realizations:")

# First bullet point - turn the new paragraph into a bulleted list item
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.InsertAfter("I have learned how to know what I don’t know, and can now think of ways find out what I don’t know. Most of the time, I’m able to control my energy to do so. But sometimes, my energy gets stuck and it’s as if I can only spend energy to learn. It leads to quite a bit of stress which and mental duress, forcing me to delve deeper into spending time simply caring for my mental health. No time for family or friends, only time for myself. It makes me need to be more selfish.")
$p.Range.ListFormat.ApplyBulletDefault()

# Second bullet point - pressing "Enter" after a list item continues the same list
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.InsertAfter("I can push myself really far. Long hours, no contact with family or friends, pure focus on how to learn and make my project succeed. But I’m getting tired of forcing myself to forget about the other important things in my life. Rather than use these things as distractions from work, I want my future work to be a distraction from them. ")

# Third bullet point - left empty, still part of the same list
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
